$d = $word.ActiveDocument

function Set-ParagraphRuns($paraIndex, $runsXml) {
    $p = $d.Paragraphs($paraIndex)
    $target = $d.Range($p.Range.Start, $p.Range.End - 1)
    $xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xmlFrag)
}

$rPrTNR = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'

# --- Rotate the three member paragraphs ---
# Paragraph 9 (currently "Petar Petrov, 2202331") -> "Aleksejs Panfilovs" + ", 2205693"
$runs9 = "<w:r>$rPrTNR<w:t>Aleksejs Panfilovs</w:t></w:r><w:r>$rPrTNR<w:t>, 2205693</w:t></w:r>"
Set-ParagraphRuns 9 $runs9

# Paragraph 10 (currently "Viktor Taskov, 2209951") -> "Petar Petrov" + ", " + "2202331"
$runs10 = "<w:r>$rPrTNR<w:t>Petar Petrov</w:t></w:r><w:r>$rPrTNR<w:t xml:space=`"preserve`">, </w:t></w:r><w:r>$rPrTNR<w:t>2202331</w:t></w:r>"
Set-ParagraphRuns 10 $runs10

# Paragraph 11 (currently "Aleksejs Panfilovs, 2205693") -> "Viktor Taskov" + ", 2209951"
$runs11 = "<w:r>$rPrTNR<w:t>Viktor Taskov</w:t></w:r><w:r>$rPrTNR<w:t>, 2209951</w:t></w:r>"
Set-ParagraphRuns 11 $runs11

Write-Host "p9:" $d.Paragraphs(9).Range.Text
Write-Host "p10:" $d.Paragraphs(10).Range.Text
Write-Host "p11:" $d.Paragraphs(11).Range.Text

# --- Move the _GoBack bookmark from the paragraph after PythonAnywhere hyperlink ---
# to the end of paragraph 9 (end of "Aleksejs Panfilovs, 2205693")
$d.Bookmarks("_GoBack").Delete()

$p9 = $d.Paragraphs(9)
$insPoint = $d.Range($p9.Range.End - 1, $p9.Range.End - 1)
$insPoint.InsertAfter("X")
$insPoint2 = $d.Range($p9.Range.End - 2, $p9.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $insPoint2)
$delRange = $d.Range($p9.Range.End - 2, $p9.Range.End - 1)
$delRange.Text = ""

$b = $d.Bookmarks("_GoBack")
Write-Host "bookmark now at:" $b.Start $b.End

# --- Delete the "Javascript Cookie" bullet paragraph ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Javascript Cookie*") {
        $d.Paragraphs($i).Range.Delete()
        break
    }
}

for ($i=1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host $i ":" $d.Paragraphs($i).Range.Text
}
